# Apply crypto price/volume updates per the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.321.62'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '2.273.58'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.23'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0957'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.91%  '

$ws.Range("E10").Value = '  +1.96%  '

$ws.Range("E11").Value = '  +3.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.70'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").Value = '2.681.52'
$ws.Range("E13").Value = '  +2.60%  '

$ws.Range("D15").Value = '54.294.70'
$ws.Range("E15").Value = '  +1.34%  '

$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").Value = '2.271.30'
$ws.Range("E17").Value = '  +1.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.18'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.42%  '

$ws.Range("E19").Value = '  +3.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '305.09'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.44'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.30%  '

$ws.Range("E22").Value = '  +0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.97'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.73%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").Value = '2.380.15'
$ws.Range("E25").Value = '  +2.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.151'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.48'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.85%  '

$ws.Range("E29").Value = '  +2.24%  '

$ws.Range("D30").Value = '0.0₃0686'
$ws.Range("E30").Value = '  +2.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.93'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.50%  '

$ws.Range("E33").Value = '  -0.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.72'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.919'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +10.22%  '

$ws.Range("E37").Value = '  +2.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.74'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.77'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.64%  '

$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("E41").Value = '  +2.72%  '

$ws.Range("E42").Value = '  +2.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.98'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.25%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0491'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +4.13%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0898'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.547'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '239.54'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.84%  '

$ws.Range("E49").Value = '  +1.13%  '

$ws.Range("E50").Value = '  +2.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.76'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.98%  '
